$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 735.7143
$ws.Range("I32").Value = 566.3333
$ws.Range("K32").Value = 566.3333
$ws.Range("M32").Value = -240.3333
$ws.Range("H55").Value = 294.75
$ws.Range("I55").Value = 273.6
$ws.Range("K55").Value = 273.6
$ws.Range("M55").Value = -59.60000000000002
$ws.Range("H86").Value = 2556.2856
$ws.Range("I86").Value = 2400.8
$ws.Range("J86").Value = 2945
$ws.Range("K86").Value = 2400.8
$ws.Range("L86").Value = 2945
$ws.Range("M86").Value = -1277.8
$ws.Range("N86").Value = -5191
$ws.Range("H89").Value = 2556.2856
$ws.Range("I89").Value = 2400.8
$ws.Range("J89").Value = 2945
$ws.Range("K89").Value = 12004
$ws.Range("L89").Value = 14725
$ws.Range("M89").Value = -6388
$ws.Range("N89").Value = -25957
$ws.Range("H132").Value = 1130.561
$ws.Range("I132").Value = 1034.6923
$ws.Range("K132").Value = 3104.0769
$ws.Range("M132").Value = -574.0769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3822.491
$ws.Range("I32").Value = 2367.9048
$ws.Range("K32").Value = 2367.9048
$ws.Range("M32").Value = -2080.9048
$ws.Range("H45").Value = 1810.1875
$ws.Range("I45").Value = 1033.4546
$ws.Range("K45").Value = 1033.4546
$ws.Range("M45").Value = -656.4546
$ws.Range("H61").Value = 3537
$ws.Range("I61").Value = 1882.7273
$ws.Range("K61").Value = 1882.7273
$ws.Range("M61").Value = -1670.7273
$ws.Range("H74").Value = 2575
$ws.Range("I74").Value = 999.25
$ws.Range("K74").Value = 999.25
$ws.Range("M74").Value = -125.25
$ws.Range("H77").Value = 2575
$ws.Range("I77").Value = 999.25
$ws.Range("K77").Value = 4996.25
$ws.Range("M77").Value = -628.25
$ws.Range("H92").Value = 43249.5
$ws.Range("J92").Value = 43249.5
$ws.Range("L92").Value = 43249.5
$ws.Range("N92").Value = -48241.5
$ws.Range("H132").Value = 1718.3182
$ws.Range("I132").Value = 1390.25
$ws.Range("K132").Value = 4170.75
$ws.Range("M132").Value = -1640.75
$ws.Range("H136").Value = 3537
$ws.Range("I136").Value = 1882.7273
$ws.Range("K136").Value = 5648.1819
$ws.Range("M136").Value = -3098.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 998.5
$ws.Range("I5").Value = 998.5
$ws.Range("K5").Value = 998.5
$ws.Range("M5").Value = -885.5
$ws.Range("H20").Value = 1999.875
$ws.Range("I20").Value = 1999.8572
$ws.Range("J20").Value = 2000
$ws.Range("K20").Value = 1999.8572
$ws.Range("L20").Value = 2000
$ws.Range("M20").Value = -1752.8572
$ws.Range("N20").Value = -2494
$ws.Range("H134").Value = 6839.5454
$ws.Range("I134").Value = 6887.7856
$ws.Range("K134").Value = 20663.3568
$ws.Range("M134").Value = -18128.3568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4999.6665
$ws.Range("I62").Value = 4999.6665
$ws.Range("K62").Value = 4999.6665
$ws.Range("M62").Value = -4375.6665
$ws.Range("H65").Value = 4999.6665
$ws.Range("I65").Value = 4999.6665
$ws.Range("K65").Value = 24998.3325
$ws.Range("M65").Value = -21878.3325
$ws.Range("H107").Value = 523.0625
$ws.Range("I107").Value = 450.9
$ws.Range("J107").Value = 643.3333
$ws.Range("K107").Value = 450.9
$ws.Range("L107").Value = 643.3333
$ws.Range("M107").Value = 1469.1
$ws.Range("N107").Value = -4483.3333
$ws.Range("H122").Value = 1611.7333
$ws.Range("I122").Value = 1463.8889
$ws.Range("K122").Value = 4391.6667
$ws.Range("M122").Value = -1941.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5064.4
$ws.Range("I80").Value = 5080.5
$ws.Range("K80").Value = 15241.5
$ws.Range("M80").Value = -14305.5
$ws.Range("H83").Value = 5064.4
$ws.Range("I83").Value = 5080.5
$ws.Range("K83").Value = 45724.5
$ws.Range("M83").Value = -41044.5
$ws.Range("H86").Value = 2247.5
$ws.Range("I86").Value = 2247.5
$ws.Range("K86").Value = 6742.5
$ws.Range("M86").Value = -5556.5
$ws.Range("H89").Value = 2247.5
$ws.Range("I89").Value = 2247.5
$ws.Range("K89").Value = 20227.5
$ws.Range("M89").Value = -14299.5
$ws.Range("H114").Value = 1951.4286
$ws.Range("J114").Value = 2696.2
$ws.Range("L114").Value = 8088.599999999999
$ws.Range("N114").Value = -14596.6
$ws.Range("H122").Value = 1064.4166
$ws.Range("I122").Value = 625.125
$ws.Range("J122").Value = 1943
$ws.Range("K122").Value = 5626.125
$ws.Range("L122").Value = 17487
$ws.Range("M122").Value = -3176.125
$ws.Range("N122").Value = -22387

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5624.5
$ws.Range("I70").Value = 6749
$ws.Range("K70").Value = 6749
$ws.Range("M70").Value = -6479
$ws.Range("H73").Value = 5624.5
$ws.Range("I73").Value = 6749
$ws.Range("K73").Value = 6749
$ws.Range("M73").Value = -5813
$ws.Range("H122").Value = 1842.7778
$ws.Range("I122").Value = 1714.7646
$ws.Range("K122").Value = 5144.293799999999
$ws.Range("M122").Value = -2694.293799999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6685
$ws.Range("J16").Value = 480.2
$ws.Range("L16").Value = 480.2
$ws.Range("N16").Value = -820.2
$ws.Range("H22").Value = 1451.9333
$ws.Range("I22").Value = 1399.6666
$ws.Range("J22").Value = 1465
$ws.Range("K22").Value = 1399.6666
$ws.Range("L22").Value = 1465
$ws.Range("M22").Value = -1104.6666
$ws.Range("N22").Value = -2055
$ws.Range("H27").Value = 1451.9333
$ws.Range("I27").Value = 1399.6666
$ws.Range("J27").Value = 1465
$ws.Range("K27").Value = 1399.6666
$ws.Range("L27").Value = 1465
$ws.Range("M27").Value = -1292.6666
$ws.Range("N27").Value = -1679
$ws.Range("H40").Value = 7862.9546
$ws.Range("I40").Value = 4989.3
$ws.Range("J40").Value = 10257.667
$ws.Range("K40").Value = 4989.3
$ws.Range("L40").Value = 10257.667
$ws.Range("M40").Value = -4853.3
$ws.Range("N40").Value = -10529.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3857.6
$ws.Range("I132").Value = 1386.6
$ws.Range("J132").Value = 6328.6
$ws.Range("K132").Value = 4159.799999999999
$ws.Range("L132").Value = 18985.8
$ws.Range("M132").Value = -1629.799999999999
$ws.Range("N132").Value = -24045.8
